$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.040.75"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.176.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.69%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.17"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.41"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.169.26"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.21"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.507"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000274"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +18.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.06"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +7.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.698.19"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.140.43"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.173.97"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.63%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.20"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.32%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "514.31"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +7.51%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.736"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +7.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.36"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.19%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.44"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.80%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.07"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +11.90%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +8.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.20"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +6.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.78"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +13.63%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Mantle"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.21"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.31"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +9.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.68"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.72"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0892"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +9.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "475.57"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.12"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0421"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.088.88"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.67"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.99%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.81%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +9.44%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +11.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.46"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0598"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +15.40%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +10.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.06"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.12%  "
